# Applies the E-column (estado) value swaps described in the commit diff.
# Rows whose "estado" value changes from "Trabaja" to "Nada"
$rowsToNada = @(6,11,24,25,61,70,71,94,95,101,121,140,141,168,173,186,187,193,213,232,233,239,262,267,278,279,285,304,309,340,360)

# Rows whose "estado" value changes from "Nada" to "Trabaja"
$rowsToTrabaja = @(12,18,19,64,65,103,132,133,172,178,179,196,224,225,241,264,270,271,288,310,316,317,338,361)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($r in $rowsToNada) {
    $ws.Cells.Item($r, 5).Value = "Nada"
}

foreach ($r in $rowsToTrabaja) {
    $ws.Cells.Item($r, 5).Value = "Trabaja"
}

$wb.Save()
